$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") table: swap the applied
#    table style from the custom "Table_0" style to the built-in style whose
#    id is {E49A603C-4CF9-4D78-9E5E-FF9A46B18C4C}.
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
for ($i = 1; $i -le $s5.Shapes.Count; $i++) {
    $shp = $s5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{E49A603C-4CF9-4D78-9E5E-FF9A46B18C4C}")
    }
}

# ---------------------------------------------------------------------------
# 2) Re-colour the deck's theme so the design swaps from the "Integral" /
#    "Red Violet" palette over to the stock "Office" palette (dk2, lt2 and
#    accent1-6 plus the hyperlink colours all change; dk1/lt1 stay
#    black/white either way).
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

$officePalette = @{
    1  = 0          # dk1      000000
    2  = 16777215   # lt1      FFFFFF
    3  = 6968388     # dk2      44546A
    4  = 15132391    # lt2      E7E6E6
    5  = 13998939    # accent1  5B9BD5
    6  = 3243501      # accent2  ED7D31
    7  = 10855845    # accent3  A5A5A5
    8  = 49407        # accent4  FFC000
    9  = 12874308    # accent5  4472C4
    10 = 4697456     # accent6  70AD47
    11 = 12673797    # hlink    0563C1
    12 = 7491477     # folHlink 954F72
}

foreach ($idx in $officePalette.Keys) {
    $colors.Item($idx).RGB = $officePalette[$idx]
}
